$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated counts for the years already present (rows 2-12) -- new search query
$ws.Range("B2").Value = 15
$ws.Range("B3").Value = 569
$ws.Range("B4").Value = 517
$ws.Range("B5").Value = 475
$ws.Range("B6").Value = 387
$ws.Range("B7").Value = 319
$ws.Range("B8").Value = 270
$ws.Range("B9").Value = 207
$ws.Range("B10").Value = 162
$ws.Range("B11").Value = 113
$ws.Range("B12").Value = 91

# New years appended below the existing data (rows 13-35)
$newData = @(
    @(2014, 93),
    @(2013, 76),
    @(2012, 51),
    @(2011, 51),
    @(2010, 33),
    @(2009, 29),
    @(2008, 28),
    @(2007, 28),
    @(2006, 13),
    @(2005, 16),
    @(2004, 15),
    @(2003, 5),
    @(2002, 14),
    @(2001, 13),
    @(2000, 11),
    @(1999, 4),
    @(1998, 5),
    @(1997, 2),
    @(1996, 6),
    @(1995, 6),
    @(1994, 3),
    @(1992, 3),
    @(1991, 1)
)

$row = 13
foreach ($pair in $newData) {
    $cellA = $ws.Cells.Item($row, 1)
    # Leading apostrophe forces the numeric-looking year to be stored as
    # text (matching the "Publication Years" column's existing text type),
    # then reset the style so no quote-prefix formatting sticks around.
    $cellA.Value = "'" + $pair[0]
    $cellA.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}
